# Contest 7 LSG vs CSK
# Row 19 on Sheet1 corresponds to Contest 7 (C19 = "LSG vs CSK").
# Fill in the points entered columns (E, H, K, N, Q, T, W) for row 19.
# The rank-lookup formulas in D, G, J, M, P, S, V (and the totals in row 35)
# recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E19").Value = 80
$ws.Range("H19").Value = 0
$ws.Range("K19").Value = 50
$ws.Range("N19").Value = 40
$ws.Range("Q19").Value = 60
$ws.Range("T19").Value = 70
$ws.Range("W19").Value = 100

$excel.CalculateFullRebuild()
